$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N, shifting the existing
# "Late" / "heading" / "Outstanding" columns one to the right.
$ws.Columns("N").Insert()

# Match the new column's width to the neighbouring "In Advance" column
# (same nominal width, but without the bestFit/autofit flag).
$ws.Columns("N").ColumnWidth = 9.83

# Make "Repayment schedule" the active sheet/tab (was "Transactions").
$ws.Activate()

# Update the active selection on the now-active sheet.
$ws.Range("S3").Select() | Out-Null
